# Update FuelPrices at 2025-04-21 09:21
# Prepend latest fetched prices (LNBSF00 / MLBSO00) to the existing date
# series: a new "today" row is added on top, the rest of the series shifts
# down by one row and keeps its original (reverse-chronological) order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New columns B/C headers (ticker symbols), matching A1's look
# ------------------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "LNBSF00"
$ws.Cells.Item(1, 3).Value = "MLBSO00"
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Column A: full reverse-chronological date series, newest first,
#    with the newest date duplicated (new fetch + existing same-day row)
# ------------------------------------------------------------------
$newDates = @(
    45764, 45764, 45763, 45762, 45761, 45758, 45757, 45756, 45755, 45754,
    45751, 45750, 45749, 45748, 45744, 45743, 45742, 45741, 45740, 45737,
    45736, 45735, 45734, 45733, 45730, 45729, 45728
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
}

# Row 28 is brand new (the series only had 26 data rows before) -> give it
# the same date-formatted style the rest of column A already carries.
$ws.Range("A2").Copy()
$ws.Range("A28").PasteSpecial(-4122)     # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3) Newly-fetched price values for the newest (top) row only
# ------------------------------------------------------------------
$ws.Cells.Item(2, 2).Value = 753.948
$ws.Cells.Item(2, 3).Value = 758.977

# All older rows have no price data yet for these new series - touch them
# so the cells exist (present, empty) without taking on any formatting.
for ($r = 3; $r -le 28; $r++) {
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Style = "Normal"
}

# ------------------------------------------------------------------
# 4) Date/time display format now includes the time component
# ------------------------------------------------------------------
$ws.Range("A2:A28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
